$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.796.96'
$ws.Range('E2').Value = '  -2.76%  '
$ws.Range('D3').Value = '2.680.20'
$ws.Range('E3').Value = '  +1.50%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.86'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.42'
$ws.Range('E6').Value = '  -6.02%  '
$ws.Range('E7').Value = '  -2.32%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.581'
$ws.Range('E9').Value = '  -2.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.39'
$ws.Range('E10').Value = '  -2.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0848'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.14'
$ws.Range('E12').Value = '  -3.44%  '
$ws.Range('D13').Value = '3.087.19'
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('E14').Value = '  +0.90%  '
$ws.Range('D15').Value = '2.677.50'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.932'
$ws.Range('E16').Value = '  -1.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '15.19'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('D18').Value = '45.862.70'
$ws.Range('E18').Value = '  -3.80%  '
$ws.Range('E19').Value = '  -1.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.90'
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.93'
$ws.Range('E21').Value = '  -2.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.75'
$ws.Range('E22').Value = '  +3.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '282.19'
$ws.Range('E23').Value = '  +4.03%  '
$ws.Range('E24').Value = '  -1.90%  '
$ws.Range('B25').Value = 'EthereumClassic'
$ws.Range('C25').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '31.13'
$ws.Range('E25').Value = '  +2.31%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.25'
$ws.Range('E26').Value = '  +1.92%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.06'
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.58'
$ws.Range('E29').Value = '  -1.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.36'
$ws.Range('E30').Value = '  -4.50%  '
$ws.Range('E31').Value = '  -6.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.24'
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('E34').Value = '  +4.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '155.96'
$ws.Range('E35').Value = '  +2.57%  '
$ws.Range('E36').Value = '  -1.51%  '
$ws.Range('E37').Value = '  -1.84%  '
$ws.Range('E38').Value = '  -2.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '25.90'
$ws.Range('E39').Value = '  +9.21%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '16.31'
$ws.Range('E41').Value = '  -3.00%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0328'
$ws.Range('E42').Value = '  -1.52%  '
$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.60'
$ws.Range('E43').Value = '  -3.65%  '
$ws.Range('E44').Value = '  -7.25%  '
$ws.Range('D45').Value = '2.154.00'
$ws.Range('E45').Value = '  -2.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '94.30'
$ws.Range('E47').Value = '  -1.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.39'
$ws.Range('E48').Value = '  -7.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '112.15'
$ws.Range('E49').Value = '  -2.31%  '
$ws.Range('D50').Value = '2.929.73'
$ws.Range('E50').Value = '  +1.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.200'
$ws.Range('E51').Value = '  -1.46%  '
